$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its values as literal text (it already
# stores things like "26.110.65" / "0.0₅8173" that are not valid numbers),
# so force text format before writing any values to avoid Excel coercing
# numeric-looking strings (e.g. "218.80") into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.110.65'
$ws.Range('E2').Value = '  -0.68%  '
$ws.Range('D3').Value = '1.656.06'
$ws.Range('E3').Value = '  -0.70%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').Value = '218.80'
$ws.Range('E5').Value = '  -0.46%  '
$ws.Range('D6').Value = '0.5308'
$ws.Range('E6').Value = '  +1.37%  '
$ws.Range('D7').Value = '1.003'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '0.2619'
$ws.Range('E8').Value = '  -2.16%  '
$ws.Range('D9').Value = '0.06343'
$ws.Range('E9').Value = '  -0.22%  '
$ws.Range('E10').Value = '  -3.43%  '
$ws.Range('D11').Value = '0.07785'
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.696.95'
$ws.Range('E12').Value = '  +1.45%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '4.497'
$ws.Range('E13').Value = '  +1.09%  '
$ws.Range('D14').Value = '0.5479'
$ws.Range('E14').Value = '  -0.53%  '
$ws.Range('D15').Value = '0.0₅8173'
$ws.Range('E15').Value = '  -1.11%  '
$ws.Range('D16').Value = '65.22'
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('D17').Value = '26.115.79'
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('E19').Value = '  -2.64%  '
$ws.Range('D20').Value = '193.16'
$ws.Range('E20').Value = '  -1.63%  '
$ws.Range('D21').Value = '10.09'
$ws.Range('E21').Value = '  -0.93%  '
$ws.Range('D22').Value = '6.026'
$ws.Range('E22').Value = '  -1.43%  '
$ws.Range('D23').Value = '1.003'
$ws.Range('E23').Value = '  -0.30%  '
$ws.Range('D24').Value = '140.24'
$ws.Range('E24').Value = '  +0.82%  '
$ws.Range('E25').Value = '  -0.25%  '
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('E27').Value = '  -0.74%  '
$ws.Range('D28').Value = '1.432'
$ws.Range('E28').Value = '  +1.18%  '
$ws.Range('D29').Value = '0.05953'
$ws.Range('E29').Value = '  -0.62%  '
$ws.Range('D30').Value = '1.282'
$ws.Range('E30').Value = '  -0.43%  '
$ws.Range('D31').Value = '3.516'
$ws.Range('E31').Value = '  -2.97%  '
$ws.Range('D32').Value = '3.242'
$ws.Range('E32').Value = '  -2.12%  '
$ws.Range('D33').Value = '1.559'
$ws.Range('E33').Value = '  -4.98%  '
$ws.Range('D34').Value = '0.9518'
$ws.Range('E34').Value = '  -3.61%  '
$ws.Range('D35').Value = '2.411'
$ws.Range('E35').Value = '  -0.55%  '
$ws.Range('D36').Value = '2.772'
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('D37').Value = '0.5667'
$ws.Range('E37').Value = '  -4.30%  '
$ws.Range('D38').Value = '0.01612'
$ws.Range('E38').Value = '  +0.73%  '
$ws.Range('D39').Value = '5.837'
$ws.Range('E39').Value = '  -3.38%  '
$ws.Range('D40').Value = '0.8463'
$ws.Range('E40').Value = '  -1.69%  '
$ws.Range('E41').Value = '  -0.17%  '
$ws.Range('D42').Value = '101.62'
$ws.Range('E42').Value = '  +1.26%  '
$ws.Range('D43').Value = '1.014.81'
$ws.Range('E43').Value = '  -1.80%  '
$ws.Range('D44').Value = '1.800.05'
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('E45').Value = '  -0.90%  '
$ws.Range('B46').Value = 'Frax'
$ws.Range('C46').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D46').Value = '1.002'
$ws.Range('E46').Value = '  -1.28%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₈104'
$ws.Range('E47').Value = '  -5.32%  '
$ws.Range('D48').Value = '0.4287'
$ws.Range('E48').Value = '  +1.38%  '
$ws.Range('D49').Value = '1.475'
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('D50').Value = '0.05158'
$ws.Range('E50').Value = '  -0.52%  '
$ws.Range('D51').Value = '7.802'
$ws.Range('E51').Value = '  -3.34%  '
